$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 3 ("Calm") so the nav colour matches row 2's nav colour (563D7C) ---
# This causes the now-unused "BDB2FF" shared string to be dropped on save.
$ws.Range("C3").Value = "563D7C"

# --- Row 4: "Green and Gold" ---
$ws.Range("C4").Value = 264653
$ws.Range("D4").Value = "2a9d8f"
$ws.Range("E4").Value = "e9c46a"
$ws.Range("F4").Value = "f4a261"
$ws.Range("G4").Value = "e76f51"
$ws.Range("A4").Value = "Green and Gold"
$ws.Range("B4").Value = $false

# --- Row 5: "Pastel" ---
$ws.Range("C5").Value = 791402
$ws.Range("D5").Value = "fec5bb"
$ws.Range("E5").Value = "fcd5ce"
$ws.Range("F5").Value = "fec89a"
$ws.Range("G5").Value = "ffd7ba"
$ws.Range("A5").Value = "Pastel"
$ws.Range("B5").Value = $false

# --- Row 6: "Ocean" ---
$ws.Range("A6").Value = "Ocean"
$ws.Range("C6").Value = "03045e"
$ws.Range("D6").Value = "023e8a"
$ws.Range("E6").Value = "0077b6"
$ws.Range("F6").Value = "48cae4"
$ws.Range("G6").Value = "90e0ef"
$ws.Range("B6").Value = $false

# --- Row 7: "Blue and Green" ---
# C7 = "184E77" typed with a forced-text leading apostrophe while the cell already
# carries a scientific number format -> text value retaining the quote-prefix flag.
# D7 = "1e6091" the same way, then we strip the quote-prefix flag back off (while
# keeping the number format) by pasting formats from a plain numFmt-only helper cell.
$ws.Range("Y1").NumberFormat = "0.00E+00"

$ws.Range("C7").NumberFormat = "0.00E+00"
$ws.Range("C7").Value = "'184E77"

$ws.Range("D7").NumberFormat = "0.00E+00"
$ws.Range("D7").Value = "'1e6091"

$ws.Range("Y1").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("Y1").Clear()

$ws.Range("E7").Value = "1a759f"
$ws.Range("F7").Value = "d9ed92"
$ws.Range("G7").Value = "b5e48c"
$ws.Range("A7").Value = "Blue and Green"
$ws.Range("B7").Value = $false

$ws.Range("B8").Select()
